$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B14").Value = "18, 45 "
$ws.Range("C14").Value = "preprocessing session 4 and trial by trial pipeline session 2, 3, 4"

$ws.Range("B12").Select()
